$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of row 5 into row 6 (same look as existing rows)
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122) # xlPasteFormats

# Now set the values for the new row 6
$ws.Range("A6").Value = 1005
$ws.Range("B6").Value = "after"
$ws.Range("C6").Value = "akashrathod290499@gmail.com"
$ws.Range("D6").Value = 97435437782

# Add hyperlink on C6
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:akashrathod290499@gmail.com")

# Update selection to D10 (matches the diff)
$ws.Range("D10").Select()
